$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.302.91'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.789.63'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.61'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5352'
$ws.Range('E7').Value = '  -0.90%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3758'
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07473'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.30'
$ws.Range('E10').Value = '  -3.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.092'
$ws.Range('E11').Value = '  -2.39%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.42'
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.086'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.238'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('D16').Value = '1.785.68'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.08'
$ws.Range('E17').Value = '  -2.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001055'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06489'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.34'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.917'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').Value = '28.326.40'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.07'
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.089'
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.32'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.24'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('D28').Value = '1.991.89'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.287'
$ws.Range('E29').Value = '  -6.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.59'
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.092'
$ws.Range('E31').Value = '  -4.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1046'
$ws.Range('E32').Value = '  +3.40%  '
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.516'
$ws.Range('E34').Value = '  -3.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2250'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06407'
$ws.Range('E36').Value = '  +2.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02273'
$ws.Range('E37').Value = '  -1.96%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.487'
$ws.Range('E39').Value = '  -4.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6146'
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.437'
$ws.Range('E41').Value = '  +3.70%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.178'
$ws.Range('E42').Value = '  +1.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.98'
$ws.Range('E43').Value = '  -4.62%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.33'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.663'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5756'
$ws.Range('E47').Value = '  -3.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.53'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.196'
$ws.Range('E49').Value = '  +4.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.927'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06845'
$ws.Range('E51').Value = '  -0.68%  '
